$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 and J1, matching style of existing header row (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# New data values in columns I and J
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9
